$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5558648
$ws.Range("I40").Value = 2079
$ws.Range("J40").Value = 8336933
$ws.Range("K40").Value = 2079
$ws.Range("L40").Value = 8336933
$ws.Range("M40").Value = -1904
$ws.Range("N40").Value = -8337283
# Row 98
$ws.Range("H98").Value = 4141.879
$ws.Range("I98").Value = 3021.3125
$ws.Range("J98").Value = 40000
$ws.Range("K98").Value = 3021.3125
$ws.Range("L98").Value = 40000
$ws.Range("M98").Value = -1523.3125
$ws.Range("N98").Value = -42996
# Row 112
$ws.Range("H112").Value = 4934.0435
$ws.Range("I112").Value = 1032.3334
$ws.Range("J112").Value = 5519.3
$ws.Range("K112").Value = 3097.0002
$ws.Range("L112").Value = 16557.9
$ws.Range("M112").Value = -1989.0002
$ws.Range("N112").Value = -18773.9
# Row 122
$ws.Range("H122").Value = 4141.879
$ws.Range("I122").Value = 3021.3125
$ws.Range("J122").Value = 40000
$ws.Range("K122").Value = 9063.9375
$ws.Range("L122").Value = 120000
$ws.Range("M122").Value = -6613.9375
$ws.Range("N122").Value = -124900
# Row 129
$ws.Range("H129").Value = 1193.8
$ws.Range("I129").Value = 965.125
$ws.Range("J129").Value = 2108.5
$ws.Range("K129").Value = 2895.375
$ws.Range("L129").Value = 6325.5
$ws.Range("M129").Value = 2104.625
$ws.Range("N129").Value = -16325.5
# Row 135
$ws.Range("H135").Value = 1111859.8
$ws.Range("I135").Value = 1111859.8
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10006738.2
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -10004203.2

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2713807.2
$ws.Range("I32").Value = 3231235.2
$ws.Range("J32").Value = 40429.418
$ws.Range("K32").Value = 3231235.2
$ws.Range("L32").Value = 40429.418
$ws.Range("M32").Value = -3230948.2
$ws.Range("N32").Value = -41003.418
# Row 37
$ws.Range("H37").Value = 10085
$ws.Range("I37").Value = 10085
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 10085
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -9812
$ws.Range("N37").ClearContents()
# Row 43
$ws.Range("H43").Value = 49999.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 49999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 49999.5
$ws.Range("N43").Value = -50625.5
$ws.Range("M43").ClearContents()
# Row 46
$ws.Range("H46").Value = 4075.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4075.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4075.5
$ws.Range("N46").Value = -4713.5
# Row 61
$ws.Range("H61").Value = 24393130
$ws.Range("I61").Value = 2105.72
$ws.Range("J61").Value = 62504108
$ws.Range("K61").Value = 2105.72
$ws.Range("L61").Value = 62504108
$ws.Range("M61").Value = -1893.72
$ws.Range("N61").Value = -62504532
# Row 102
$ws.Range("H102").Value = 3344.6667
$ws.Range("I102").Value = 2827.4443
$ws.Range("J102").Value = 7999.6665
$ws.Range("K102").Value = 2827.4443
$ws.Range("L102").Value = 7999.6665
$ws.Range("M102").Value = -1205.4443
$ws.Range("N102").Value = -11243.6665
# Row 109
$ws.Range("H109").Value = 59340.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 59340.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 59340.5
$ws.Range("N109").Value = -62114.5
# Row 112
$ws.Range("H112").Value = 59279
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 59279
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 59279
$ws.Range("N112").Value = -62233
# Row 119
$ws.Range("H119").Value = 51979
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 51979
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 51979
$ws.Range("N119").Value = -61655
# Row 122
$ws.Range("H122").Value = 1965.0303
$ws.Range("I122").Value = 1465.6552
$ws.Range("J122").Value = 5585.5
$ws.Range("K122").Value = 4396.9656
$ws.Range("L122").Value = 16756.5
$ws.Range("M122").Value = -1946.9656
$ws.Range("N122").Value = -21656.5
# Row 135
$ws.Range("H135").Value = 1411793.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1411793.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 1411793.5
$ws.Range("N135").Value = -1421933.5
# Row 136
$ws.Range("H136").Value = 24393130
$ws.Range("I136").Value = 2105.72
$ws.Range("J136").Value = 62504108
$ws.Range("K136").Value = 6317.16
$ws.Range("L136").Value = 187512324
$ws.Range("M136").Value = -3767.16
$ws.Range("N136").Value = -187517424

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6949929
$ws.Range("I20").Value = 11115874
$ws.Range("J20").Value = 6687.1113
$ws.Range("K20").Value = 11115874
$ws.Range("L20").Value = 6687.1113
$ws.Range("M20").Value = -11115627
$ws.Range("N20").Value = -7181.1113
# Row 57
$ws.Range("H57").Value = 69973.164
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 69973.164
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 69973.164
$ws.Range("N57").Value = -71413.164
# Row 132
$ws.Range("H132").Value = 107500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 107500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 107500
$ws.Range("N132").Value = -117620
# Row 134
$ws.Range("H134").Value = 5438022.5
$ws.Range("I134").Value = 8621520
$ws.Range("J134").Value = 7349.5884
$ws.Range("K134").Value = 25864560
$ws.Range("L134").Value = 22048.7652
$ws.Range("M134").Value = -25862025
$ws.Range("N134").Value = -27118.7652
# Row 136
$ws.Range("H136").Value = 69973.164
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 69973.164
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 69973.164
$ws.Range("N136").Value = -80173.164

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2
$ws.Range("N4").Value = -226
# Row 31
$ws.Range("H31").Value = 6117.55
$ws.Range("I31").Value = 2730.375
$ws.Range("J31").Value = 6964.3438
$ws.Range("K31").Value = 2730.375
$ws.Range("L31").Value = 6964.3438
$ws.Range("M31").Value = -2435.375
$ws.Range("N31").Value = -7554.3438
# Row 34
$ws.Range("H34").Value = 6117.55
$ws.Range("I34").Value = 2730.375
$ws.Range("J34").Value = 6964.3438
$ws.Range("K34").Value = 2730.375
$ws.Range("L34").Value = 6964.3438
$ws.Range("M34").Value = -2528.375
$ws.Range("N34").Value = -7368.3438
# Row 42
$ws.Range("H42").Value = 49999.5
$ws.Range("I42").Value = 49999
$ws.Range("J42").Value = 50000
$ws.Range("K42").Value = 49999
$ws.Range("L42").Value = 50000
$ws.Range("M42").Value = -49406
$ws.Range("N42").Value = -51186
# Row 98
$ws.Range("H98").Value = 42185.715
$ws.Range("I98").Value = 40000
$ws.Range("J98").Value = 55300
$ws.Range("K98").Value = 40000
$ws.Range("L98").Value = 55300
$ws.Range("M98").Value = -37754
$ws.Range("N98").Value = -59792
# Row 122
$ws.Range("H122").Value = 4375.909
$ws.Range("I122").Value = 3930.6667
$ws.Range("J122").Value = 4542.875
$ws.Range("K122").Value = 11792.0001
$ws.Range("L122").Value = 13628.625
$ws.Range("M122").Value = -9342.000100000001
$ws.Range("N122").Value = -18528.625
# Row 123
$ws.Range("H123").Value = 78000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 78000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 78000
$ws.Range("N123").Value = -87800
# Row 141
$ws.Range("H141").Value = 203927
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 203927
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 203927
$ws.Range("N141").Value = -214287

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 96131.38
$ws.Range("I2").Value = 142.6
$ws.Range("J2").Value = 183393.9
$ws.Range("K2").Value = 855.5999999999999
$ws.Range("L2").Value = 1100363.4
$ws.Range("M2").Value = -742.5999999999999
$ws.Range("N2").Value = -1100589.4
# Row 4
$ws.Range("H4").Value = 50321190
$ws.Range("I4").Value = 62527910
$ws.Range("J4").Value = 3935661.5
$ws.Range("K4").Value = 187583730
$ws.Range("L4").Value = 11806984.5
$ws.Range("M4").Value = -187583618
$ws.Range("N4").Value = -11807208.5
# Row 38
$ws.Range("H38").Value = 55555584
$ws.Range("I38").Value = 24.75
$ws.Range("J38").Value = 100000030
$ws.Range("K38").Value = 74.25
$ws.Range("L38").Value = 300000090
$ws.Range("M38").Value = 272.75
$ws.Range("N38").Value = -300000784
# Row 55
$ws.Range("H55").Value = 8348900
$ws.Range("I55").Value = 9602.666999999999
$ws.Range("J55").Value = 11128666
$ws.Range("K55").Value = 28808.001
$ws.Range("L55").Value = 33385998
$ws.Range("M55").Value = -28631.001
$ws.Range("N55").Value = -33386352
# Row 64
$ws.Range("H64").Value = 2276
$ws.Range("I64").Value = 1128.3334
$ws.Range("J64").Value = 3997.5
$ws.Range("K64").Value = 3385.0002
$ws.Range("L64").Value = 11992.5
$ws.Range("M64").Value = -3115.0002
$ws.Range("N64").Value = -12532.5
# Row 67
$ws.Range("H67").Value = 2276
$ws.Range("I67").Value = 1128.3334
$ws.Range("J67").Value = 3997.5
$ws.Range("K67").Value = 3385.0002
$ws.Range("L67").Value = 11992.5
$ws.Range("M67").Value = -2449.0002
$ws.Range("N67").Value = -13864.5
# Row 113
$ws.Range("H113").Value = 5021.3335
$ws.Range("I113").Value = 1899
$ws.Range("J113").Value = 5913.4287
$ws.Range("K113").Value = 5697
$ws.Range("L113").Value = 17740.2861
$ws.Range("M113").Value = -3527
$ws.Range("N113").Value = -22080.2861
# Row 140
$ws.Range("H140").Value = 104592.9
$ws.Range("I140").Value = 201955.6
$ws.Range("J140").Value = 7230.2
$ws.Range("K140").Value = 605866.8
$ws.Range("L140").Value = 21690.6
$ws.Range("M140").Value = -600686.8
$ws.Range("N140").Value = -32050.6

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 121
$ws.Range("H121").Value = 56134
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 56134
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 56134
$ws.Range("N121").Value = -59628
# Row 122
$ws.Range("H122").Value = 3460778.2
$ws.Range("I122").Value = 4274355.5
$ws.Range("J122").Value = 3075
$ws.Range("K122").Value = 12823066.5
$ws.Range("L122").Value = 9225
$ws.Range("M122").Value = -12820616.5
$ws.Range("N122").Value = -14125
# Row 132
$ws.Range("H132").Value = 4296.6665
$ws.Range("I132").Value = 3247.5
$ws.Range("J132").Value = 4821.25
$ws.Range("K132").Value = 9742.5
$ws.Range("L132").Value = 14463.75
$ws.Range("M132").Value = -7212.5
$ws.Range("N132").Value = -19523.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 90000
$ws.Range("I2").Value = 80000
$ws.Range("J2").Value = 100000
$ws.Range("K2").Value = 80000
$ws.Range("L2").Value = 100000
$ws.Range("M2").Value = -79888
$ws.Range("N2").Value = -100224
# Row 16
$ws.Range("H16").Value = 1076.45
$ws.Range("I16").Value = 973.94446
$ws.Range("J16").Value = 1999
$ws.Range("K16").Value = 973.94446
$ws.Range("L16").Value = 1999
$ws.Range("M16").Value = -803.94446
$ws.Range("N16").Value = -2339
# Row 110
$ws.Range("H110").Value = 250041070
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 250041070
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 250041070
$ws.Range("N110").Value = -250049250
# Row 119
$ws.Range("H119").Value = 56134
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 56134
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 56134
$ws.Range("N119").Value = -65810
# Row 132
$ws.Range("H132").Value = 10874272
$ws.Range("I132").Value = 18521682
$ws.Range("J132").Value = 6897.737
$ws.Range("K132").Value = 55565046
$ws.Range("L132").Value = 20693.211
$ws.Range("M132").Value = -55562516
$ws.Range("N132").Value = -25753.211

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 119
$ws.Range("H119").Value = 56133.5
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 56133.5
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 56133.5
$ws.Range("N119").Value = -65809.5
# Row 122
$ws.Range("H122").Value = 114801.414
$ws.Range("I122").Value = 192647.72
$ws.Range("J122").Value = 5816.6
$ws.Range("K122").Value = 577943.16
$ws.Range("L122").Value = 17449.8
$ws.Range("M122").Value = -575493.16
$ws.Range("N122").Value = -22349.8
# Row 132
$ws.Range("H132").Value = 5318.2104
$ws.Range("I132").Value = 5157.4443
$ws.Range("J132").Value = 5712.8184
$ws.Range("K132").Value = 15472.3329
$ws.Range("L132").Value = 17138.4552
$ws.Range("M132").Value = -12942.3329
$ws.Range("N132").Value = -22198.4552
# Row 136
$ws.Range("H136").Value = 17416968
$ws.Range("I136").Value = 23810822
$ws.Range("J136").Value = 633100.9399999999
$ws.Range("K136").Value = 71432466
$ws.Range("L136").Value = 1899302.82
$ws.Range("M136").Value = -71429916
$ws.Range("N136").Value = -1904402.82
